# Update the Agenda worksheet (Sheet1) with the new schedule / intro deck timings,
# per commit "updated agenda and intro power point".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 ---
$ws.Range("A2").Value2 = "8am-8:45am"
$ws.Range("B2").Value2 = "Intro"
$ws.Range("C2").Value2 = "Ryan"

# --- Row 3 ---
$ws.Range("A3").Value2 = "8:45am-9:45am"
$ws.Range("B3").Value2 = "Code Smells"
$ws.Range("C3").Value2 = "Patrick"

# --- Row 4 ---
$ws.Range("A4").Value2 = "9:45am-10am"
$ws.Range("B4").Value2 = "Break"
$ws.Range("C4").Value2 = "All"

# --- Row 5 ---
$ws.Range("A5").Value2 = "10:00am-12pm"
$ws.Range("B5").Value2 = "Code Challenges"
$ws.Range("C5").Value2 = "All"

# --- Row 6 ---
$ws.Range("A6").Value2 = "12pm-1pm"
$ws.Range("B6").Value2 = "Lunch"

# --- Row 7 ---
$ws.Range("A7").Value2 = "1pm-1:45pm"
$ws.Range("B7").Value2 = "Writing Clean Code"
$ws.Range("C7").Value2 = "Cory"

# --- Row 8 ---
$ws.Range("A8").Value2 = "1:45pm-2:45pm"
$ws.Range("B8").Value2 = "Refactoring practice and assistance"
$ws.Range("C8").Value2 = "All"

# --- Row 9 ---
$ws.Range("A9").Value2 = "2:45pm-3pm"
$ws.Range("B9").Value2 = "Break"
$ws.Range("C9").Value2 = "All"

# --- Row 10 ---
$ws.Range("A10").Value2 = "3pm-4pm"
$ws.Range("B10").Value2 = "Attendee code sharing"
$ws.Range("C10").Value2 = "All"

# --- Row 11 ---
$ws.Range("A11").Value2 = "4pm-5pm"
$ws.Range("B11").Value2 = "Conclusions"
$ws.Range("C11").Value2 = "Patrick"

# Row 5 used to have an explicit (taller) custom height for the old, longer
# text; the new copy fits on one line, so drop back to the default/auto row height.
$ws.Rows.Item(5).AutoFit()

# Column A needs to widen slightly to fit the longest new time range text.
$ws.Columns.Item(1).ColumnWidth = 13.83

# Update the selection left over from editing, to span the refreshed table body.
$ws.Range("A2:B11").Select() | Out-Null
